# Update feed logs and data lake files
# Append two new log rows to the feed logs sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(28, 3).Value = "2024-06-15 11:10:40"
$ws.Cells.Item(28, 4).Value = 200
$ws.Cells.Item(28, 5).Value = 5

$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = 2
$ws.Cells.Item(29, 3).Value = "2024-06-15 11:10:40"
$ws.Cells.Item(29, 4).Value = 200
$ws.Cells.Item(29, 5).Value = 0
